$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B/C (text) and E (percentage strings) columns: plain text assignment is safe.
# D (price) column: force Text format so numeric-looking strings are not
# reinterpreted as numbers/dates by Excel, matching the source inlineStr cells.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '65.160.88'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +2.69%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.482.09'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '579.71'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.04%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '162.69'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +4.36%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.614'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +12.62%  '
$ws.Range("E8").Value = '  +0.03%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '3.483.56'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +2.10%  '
$ws.Range("E10").Value = '  -1.84%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.125'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +2.90%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.448'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +3.78%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.085.95'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +2.14%  '
$ws.Range("E14").Value = '  +0.43%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.0000193'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.13%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '28.61'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +5.49%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '65.228.33'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +2.66%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.488.36'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +2.83%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.48'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +3.68%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.42'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.37%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '382.94'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.85%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '8.20'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.95%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.553'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +4.73%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '72.75'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +1.42%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.09%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.0000120'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("E27").Value = '  +7.99%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.179'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '1.54'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +12.63%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '6.20'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.59%  '
$ws.Range("E32").Value = '  +2.81%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '23.71'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.99%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '7.20'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +13.72%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '161.98'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.35%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.93'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +6.68%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0782'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +2.78%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.006.78'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.10%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '26.84'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.84%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '6.75'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +5.92%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '4.58'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +6.21%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.0323'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +2.33%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '42.90'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +2.43%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.782'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.46%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '26.04'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +12.00%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.11'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +4.24%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '318.55'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +8.15%  '
$ws.Range("E49").Value = '  +6.95%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.876'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +5.32%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.67'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +5.13%  '
